$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark target cells as Text so numeric-looking strings (e.g. "19.61", "26.718.69")
# are preserved verbatim instead of being parsed into numbers.
$targetCells = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "E6", "E7", "E8", "E9", "D10", "E10", "E11", "D12", "E12", "D13", "E13", "E15", "D16", "E16", "D17", "E17", "E18", "D19", "D20", "E20", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "E25", "E26", "E27", "D28", "E28", "D29", "E29", "D30", "E30", "D31", "E31", "E32", "D33", "E33", "B34", "C34", "D34", "E34", "B35", "C35", "D35", "E35", "E36", "D37", "E37", "E38", "E39", "E40", "E41", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "E46", "E47", "D48", "E48", "D49", "E49", "E51")
foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.718.69"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "1.598.46"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "211.72"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").Value = "19.61"
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("D12").Value = "1.822.85"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "1.588.21"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("D16").Value = "65.02"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "0.0₃0739"
$ws.Range("E17").Value = "  -2.46%  "
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").Value = "208.59"
$ws.Range("D20").Value = "7.14"
$ws.Range("E20").Value = "  +1.30%  "
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("D22").Value = "2.24"
$ws.Range("E22").Value = "  -3.60%  "
$ws.Range("D23").Value = "9.03"
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("D24").Value = "144.09"
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("D28").Value = "15.34"
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").Value = "0.0507"
$ws.Range("E29").Value = "  -2.09%  "
$ws.Range("D30").Value = "1.16"
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("D31").Value = "3.26"
$ws.Range("E31").Value = "  +0.96%  "
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("D33").Value = "1.275.77"
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "1.25"
$ws.Range("E34").Value = "  +16.80%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.48"
$ws.Range("E35").Value = "  +1.53%  "
$ws.Range("E36").Value = "  +0.59%  "
$ws.Range("D37").Value = "0.588"
$ws.Range("E37").Value = "  -4.62%  "
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("E40").Value = "  +0.57%  "
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("D43").Value = "62.54"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").Value = "1.734.65"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").Value = "90.48"
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("E46").Value = "  +0.84%  "
$ws.Range("E47").Value = "  +1.51%  "
$ws.Range("D48").Value = "0.0512"
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("D49").Value = "7.55"
$ws.Range("E49").Value = "  +3.34%  "
$ws.Range("E51").Value = "  +1.62%  "

# Restore default (General) cell style now that the text values are stored,
# so no residual text-format styling is left on the cells.
foreach ($addr in $targetCells) {
    $ws.Range($addr).Style = "Normal"
}
